$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 3.182878228561681, 86.29678392075563, 2938.103010863317, 6.48142807727062, 1, 3034.064101089905),
    @(3, 1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 4.371470058157054),
    @(4, 1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 3.811642989160245),
    @(5, 0.1554434735375247, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 3.021299490525381),
    @(6, 0.7287194209349384, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 3.594575437922795),
    @(7, 3.182878228561681, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1, 4.733082622659194),
    @(8, 3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
